# Update the "model selection" sheet with new model-comparison metrics
# (re-run using 5-min data; model update).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model selection")

# --- XGBoost Regressor tuned (row 6) ---
$ws.Range("E6").Value = 0.07
$ws.Range("F6").Value = 0.1
$ws.Range("G6").Value = 0.92
$ws.Range("H6").Value = 0.84

# --- XGBoost Regressor default (row 7) ---
$ws.Range("E7").Value = 0.08
$ws.Range("F7").Value = 0.11
$ws.Range("G7").Value = 0.91
$ws.Range("H7").Value = 0.82

# --- Linear Regression (row 8) : MAE cleared ---
$ws.Range("E8").ClearContents()

# --- Random Forest Regressor (row 11) ---
$ws.Range("E11").Value = 0.07
$ws.Range("F11").Value = 0.11
$ws.Range("G11").Value = 0.92
$ws.Range("H11").Value = 0.83

# --- Support Vector Regressor (row 13) : MAE cleared ---
$ws.Range("E13").ClearContents()

# --- Neural Network Regressor (row 14) ---
$ws.Range("E14").Value = 0.09
$ws.Range("F14").Value = 0.13
$ws.Range("G14").Value = 0.87
$ws.Range("H14").Value = 0.73

# Update the active selection on the sheet to match the latest edit location
$ws.Activate()
$ws.Range("H9").Select()
